$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one
# day (46075 -> 46076) for every data row (rows 2 through 149).
for ($r = 2; $r -le 149; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value = 46076
    }
}
